# Update Name of Algo - apply updated KNN imputation results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.601999999999999
$ws.Range("D3").Value = -7.287999999999999
$ws.Range("D5").Value = -8.02
$ws.Range("B9").Value = 6.795
$ws.Range("D11").Value = -8.15
$ws.Range("D12").Value = -8.088999999999999
$ws.Range("B13").Value = 6.472
$ws.Range("B16").Value = 5.787
$ws.Range("B18").Value = 6.313
$ws.Range("B20").Value = 6.661
$ws.Range("D21").Value = -7.895999999999999
